$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Notes") gets an "Act" (Acceptance Criteria) marker added for every
# story row that doesn't already carry notes text. Row 3 (US9) already has its
# own Acceptance Criteria note in D3, so it is left untouched.
$rows = @(2, 4, 5, 6, 7, 8)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "Act"
}
